# Applies the "Normalize preset option counts and add Parry Only profile" edit.
#
# Summary of changes:
#  - Menu Mock sheet (col E) and Providers sheet (col B): update several
#    pipe-delimited preset option lists (add "Very Rare", "Very Short",
#    remove "Cinematic", add "Parry Only", add "12.5%"/"37.5%" chance steps).
#  - Menu Mock sheet (col D): normalize several per-trigger Duration/
#    Cooldown/Smoothing default values to new numbers.

$wb = $excel.ActiveWorkbook
$wsMenu = $wb.Worksheets.Item("Menu Mock")
$wsProv = $wb.Worksheets.Item("Providers")

# --- Menu Mock sheet, column E (preset option lists) ---------------------
$menuEChanges = @{
    10 = "Off (Cooldown Only) | Very Rare | Rare | Standard | Frequent"
    12 = "Very Short | Short | Standard | Long | Extended"
    13 = "Very Snappy | Snappy | Standard | Smooth | Ultra Smooth"
    14 = "All Triggers | Kills Only | Highlights | Last Enemy Only | Parry Only"
    40 = "12.5% | 15% | 25% | 30% | 35% | 36% | 37.5% | 45% | 50% | 54% | 60% | 70% | 75% | 84% | 90% | 100%"
    48 = "12.5% | 15% | 25% | 30% | 35% | 36% | 37.5% | 45% | 50% | 54% | 60% | 70% | 75% | 84% | 90% | 100%"
    56 = "12.5% | 15% | 25% | 30% | 35% | 36% | 37.5% | 45% | 50% | 54% | 60% | 70% | 75% | 84% | 90% | 100%"
    64 = "12.5% | 15% | 25% | 30% | 35% | 36% | 37.5% | 45% | 50% | 54% | 60% | 70% | 75% | 84% | 90% | 100%"
    72 = "12.5% | 15% | 25% | 30% | 35% | 36% | 37.5% | 45% | 50% | 54% | 60% | 70% | 75% | 84% | 90% | 100%"
    86 = "12.5% | 15% | 25% | 30% | 35% | 36% | 37.5% | 45% | 50% | 54% | 60% | 70% | 75% | 84% | 90% | 100%"
}

foreach ($row in $menuEChanges.Keys) {
    $wsMenu.Range("E$row").Value = $menuEChanges[$row]
}

# --- Menu Mock sheet, column D (normalized default values) ---------------
$menuDChanges = @{
    42 = "1.0s"
    43 = "5.0s"
    44 = "8x"
    50 = "1.5s"
    51 = "5.0s"
    52 = "8x"
    58 = "1.5s"
    59 = "5.0s"
    60 = "8x"
    66 = "2.0s"
    67 = "4.0s"
    68 = "6x"
    74 = "3.0s"
    76 = "4x"
    81 = "5.0s"
    82 = "45.0s"
    83 = "4x"
    88 = "1.2s"
    89 = "7.0s"
    90 = "10x"
}

foreach ($row in $menuDChanges.Keys) {
    $wsMenu.Range("D$row").Value = $menuDChanges[$row]
}

# --- Providers sheet, column B (preset option lists, mirrors column E) ---
$provBChanges = @{
    5  = "Off (Cooldown Only) | Very Rare | Rare | Standard | Frequent"
    9  = "12.5% | 15% | 25% | 30% | 35% | 36% | 37.5% | 45% | 50% | 54% | 60% | 70% | 75% | 84% | 90% | 100%"
    15 = "Very Short | Short | Standard | Long | Extended"
    27 = "Very Snappy | Snappy | Standard | Smooth | Ultra Smooth"
    30 = "All Triggers | Kills Only | Highlights | Last Enemy Only | Parry Only"
}

foreach ($row in $provBChanges.Keys) {
    $wsProv.Range("B$row").Value = $provBChanges[$row]
}
